# Append a new row (row 29) to the historical log with June 28th, 2020 data.
# Raw and Clean Data from SSA for June 28th

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: date label as literal text "2020-06-28" -----------------
# A direct Range.Value assignment of a date-like string gets auto-converted
# by Excel into a date serial number (and picks up a quote-prefixed style).
# To keep it as plain text (matching the rest of column A, which stores
# dates as shared strings with no special style), enter it as a text
# formula and then convert it to a literal value via copy / paste-values.
$ws.Range("A29").Formula = "=""2020-06-28"""
$ws.Range("A29").Copy()
$ws.Range("A29").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# --- Columns B-F: raw counts and the computed percentage ---------------
$ws.Range("B29").Value = 216852
$ws.Range("C29").Value = 275203
$ws.Range("D29").Value = 64143
$ws.Range("E29").Value = 26648
$ws.Range("F29").Value = 31.11
